# SectorCOFOG.xlsx — rotate the last three populated columns (E, F, G) on
# every row (including the header row):
#   new E = old F
#   new F = old G   (cleared if old G was empty)
#   new G = old E
#
# This matches the upstream codeforIATI/codelists regeneration where the
# "group"/"class"/"category-code" columns were re-emitted in a different
# order (class, category-code, group) while the underlying values per row
# stayed the same set, just shifted across the three columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)
    $gCell = $ws.Cells.Item($r, 7)

    $eVal = $eCell.Text
    $fVal = $fCell.Text
    $gVal = $gCell.Text

    $hasG = $gVal -ne ""

    # Write order doesn't matter since we captured all three originals
    # already, but assign from "last" column backwards to be safe anyway.
    $eCell.Value = $fVal
    if ($hasG) {
        $fCell.Value = $gVal
    } else {
        $fCell.Value = ""
    }
    $gCell.Value = $eVal
}
